$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 (shifts existing rows 78..163 down to 79..164)
$ws.Range("A78").EntireRow.Insert()

# Populate the new row 78 with its data
$ws.Range("A78").Value = 4
$ws.Range("B78").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C78").Value = "Los Lagos"
$ws.Range("D78").Value = 44494
$ws.Range("E78").Value = 10
$ws.Range("F78").Value = 100112043
$ws.Range("G78").Value = "Pepino ensalada"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 200
$ws.Range("K78").Value = 13000
$ws.Range("L78").Value = 13000
$ws.Range("M78").Value = 13000
$ws.Range("N78").Value = "$/caja 60 unidades"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 217
$ws.Range("Q78").Value = 60
$ws.Range("R78").Value = "Hortaliza"
